$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.731.54'
$ws.Range("E2").Value = '  +2.08%  '
$ws.Range("D3").Value = '3.970.89'
$ws.Range("E3").Value = '  +0.71%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.36%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '532.68'
$ws.Range("E5").Value = '  +7.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.11'
$ws.Range("E6").Value = '  -1.25%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.619'
$ws.Range("E7").Value = '  -0.88%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.998'
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.735'
$ws.Range("E9").Value = '  +0.25%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.174'
$ws.Range("E10").Value = '  -1.23%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000344'
$ws.Range("E11").Value = '  -1.91%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '42.69'
$ws.Range("E12").Value = '  -1.60%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.52'
$ws.Range("E13").Value = '  +0.50%  '
$ws.Range("D14").Value = '4.577.47'
$ws.Range("E14").Value = '  +0.02%  '
$ws.Range("D15").Value = '4.017.44'
$ws.Range("E15").Value = '  +1.45%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.44'
$ws.Range("E16").Value = '  +7.46%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.25'
$ws.Range("E17").Value = '  +0.14%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.22'
$ws.Range("E18").Value = '  +2.33%  '
$ws.Range("E19").Value = '  -1.91%  '
$ws.Range("D20").Value = '70.585.22'
$ws.Range("E20").Value = '  +1.76%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '434.91'
$ws.Range("E21").Value = '  -0.74%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.54'
$ws.Range("E22").Value = '  +2.77%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '89.41'
$ws.Range("E23").Value = '  +0.58%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '14.28'
$ws.Range("E24").Value = '  -2.44%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.08'
$ws.Range("E25").Value = '  +7.02%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.82'
$ws.Range("E26").Value = '  -1.74%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.72'
$ws.Range("E27").Value = '  -3.74%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '36.79'
$ws.Range("E28").Value = '  -1.05%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '695.03'
$ws.Range("E29").Value = '  -0.72%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '13.49'
$ws.Range("E30").Value = '  +0.67%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.127'
$ws.Range("E31").Value = '  -2.36%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.89'
$ws.Range("E32").Value = '  +0.96%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.90'
$ws.Range("E33").Value = '  +13.35%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '66.98'
$ws.Range("E34").Value = '  +7.07%  '
$ws.Range("D35").Value = '0.0₃0907'
$ws.Range("E35").Value = '  +1.71%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.448'
$ws.Range("E36").Value = '  -2.64%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '40.11'
$ws.Range("E37").Value = '  -2.38%  '
$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.150'
$ws.Range("E38").Value = '  -0.38%  '
$ws.Range("B39").Value = 'ThetaToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.44'
$ws.Range("E39").Value = '  +11.19%  '
$ws.Range("E40").Value = '  +0.36%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("E41").Value = '  -0.30%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0484'
$ws.Range("E42").Value = '  -0.88%  '
$ws.Range("B43").Value = 'Fetch.AI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.87'
$ws.Range("E43").Value = '  -1.39%  '
$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.16'
$ws.Range("E44").Value = '  +4.26%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.55'
$ws.Range("E45").Value = '  +5.16%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.26'
$ws.Range("E46").Value = '  +8.97%  '
$ws.Range("E47").Value = '  -0.17%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.000290'
$ws.Range("E48").Value = '  +22.24%  '
$ws.Range("B49").Value = 'THORChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.13'
$ws.Range("E49").Value = '  +4.17%  '
$ws.Range("B50").Value = 'LidoDAOToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.37'
$ws.Range("E50").Value = '  -0.54%  '
$ws.Range("D51").Value = '0.0₆0349'
$ws.Range("E51").Value = '  -0.37%  '
